$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move "**basic" label from C24 down to C25
$ws.Range("C24").Value = $null
$ws.Range("C25").Value = "**basic"

# Add the new "advanced" row of results
$ws.Range("D26").Value = 310.75
$ws.Range("E26").Value = 10
$ws.Range("F26").Value = 16
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = 8
$ws.Range("I26").Value = "960 by 720"
$ws.Range("C26").Value = "**advanced"

# Update the view state (top-left visible cell + active selection)
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("C23").Select()
